$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, [string]$val) {
    # Force text storage so numeric-looking strings (e.g. "1.000")
    # are not auto-converted to numbers by Excel, matching the
    # original inlineStr/text cells, then restore default styling.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "26.533.09"
$ws.Range("E2").Value = "  +2.34%  "
Set-TextValue $ws.Range("D3") "1.681.94"
$ws.Range("E3").Value = "  +3.19%  "
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  +0.00%  "
Set-TextValue $ws.Range("D5") "216.78"
$ws.Range("E5").Value = "  +2.48%  "
$ws.Range("E6").Value = "  +1.67%  "
Set-TextValue $ws.Range("D7") "1.000"
$ws.Range("E7").Value = "  +0.07%  "
Set-TextValue $ws.Range("D8") "0.2683"
$ws.Range("E8").Value = "  +3.89%  "
Set-TextValue $ws.Range("D9") "0.06400"
$ws.Range("E9").Value = "  +1.90%  "
Set-TextValue $ws.Range("D10") "21.67"
$ws.Range("E10").Value = "  +5.49%  "
Set-TextValue $ws.Range("D11") "0.07798"
$ws.Range("E11").Value = "  +2.99%  "
Set-TextValue $ws.Range("D12") "1.683.26"
$ws.Range("E12").Value = "  +3.34%  "
Set-TextValue $ws.Range("D13") "4.504"
$ws.Range("E13").Value = "  +1.85%  "
Set-TextValue $ws.Range("D14") "0.5577"
$ws.Range("E14").Value = "  +0.97%  "
Set-TextValue $ws.Range("D15") "0.0₅8348"
$ws.Range("E15").Value = "  +3.96%  "
Set-TextValue $ws.Range("D16") "65.78"
$ws.Range("E16").Value = "  +1.49%  "
Set-TextValue $ws.Range("D17") "26.545.04"
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("E18").Value = "  +0.06%  "
Set-TextValue $ws.Range("D19") "4.772"
$ws.Range("E19").Value = "  +1.99%  "
Set-TextValue $ws.Range("D20") "195.16"
$ws.Range("E20").Value = "  +5.05%  "
$ws.Range("E21").Value = "  +2.19%  "
Set-TextValue $ws.Range("D22") "6.349"
$ws.Range("E22").Value = "  +3.65%  "
$ws.Range("E23").Value = "  +0.17%  "
Set-TextValue $ws.Range("D24") "143.34"
$ws.Range("E24").Value = "  -1.23%  "
Set-TextValue $ws.Range("D25") "0.1286"
$ws.Range("E25").Value = "  +5.64%  "
Set-TextValue $ws.Range("D26") "7.438"
$ws.Range("E26").Value = "  +0.71%  "
Set-TextValue $ws.Range("D27") "16.38"
$ws.Range("E27").Value = "  +4.30%  "
Set-TextValue $ws.Range("D28") "1.430"
$ws.Range("E28").Value = "  +4.92%  "
Set-TextValue $ws.Range("D29") "0.06176"
$ws.Range("E29").Value = "  +4.72%  "
$ws.Range("E30").Value = "  +2.21%  "
Set-TextValue $ws.Range("D31") "3.610"
$ws.Range("E31").Value = "  +5.54%  "
Set-TextValue $ws.Range("D32") "3.460"
$ws.Range("E32").Value = "  +2.04%  "
$ws.Range("E33").Value = "  +4.01%  "
Set-TextValue $ws.Range("D34") "1.009"
$ws.Range("E34").Value = "  +2.96%  "
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("E36").Value = "  +2.31%  "
Set-TextValue $ws.Range("D37") "0.5733"
$ws.Range("E37").Value = "  -0.98%  "
Set-TextValue $ws.Range("D38") "0.01637"
$ws.Range("E38").Value = "  +2.18%  "
Set-TextValue $ws.Range("D39") "6.030"
$ws.Range("E39").Value = "  +6.36%  "
Set-TextValue $ws.Range("D40") "1.077.03"
$ws.Range("E40").Value = "  +3.94%  "
Set-TextValue $ws.Range("D41") "0.8604"
$ws.Range("E41").Value = "  +1.43%  "
$ws.Range("E42").Value = "  -0.25%  "
Set-TextValue $ws.Range("D43") "100.12"
$ws.Range("E43").Value = "  +0.25%  "
Set-TextValue $ws.Range("D44") "1.823.87"
$ws.Range("E44").Value = "  +2.90%  "
Set-TextValue $ws.Range("D45") "0.0₈111"
$ws.Range("E45").Value = "  +4.12%  "
Set-TextValue $ws.Range("D46") "57.03"
$ws.Range("E46").Value = "  +3.48%  "
Set-TextValue $ws.Range("D47") "8.147"
$ws.Range("E47").Value = "  +1.37%  "
Set-TextValue $ws.Range("D48") "1.003"
$ws.Range("E48").Value = "  -0.34%  "
Set-TextValue $ws.Range("D49") "0.05208"
$ws.Range("E49").Value = "  +0.84%  "

# Row 50: Aptos -> RenderToken (new coin inserted)
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D50") "1.471"
$ws.Range("E50").Value = "  +5.96%  "

# Row 51: Mantle -> Aptos (Mantle dropped off the list, Aptos shifts down)
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D51") "6.030"
$ws.Range("E51").Value = "  +2.88%  "
